# moved combat into its own module. changed combat formulas.
#
# Adds a new "Sheet2" (combat-log module) after "Sheet1", populates it with
# the combat log / stat rows, and makes it the active sheet (mirrors the
# author moving the combat-log section out of Sheet1 into its own sheet).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Insert the new sheet right after Sheet1 -> becomes "Sheet2"
$ws2 = $wb.Worksheets.Add($null, $ws1)

# Combat log lines
$ws2.Range("B4").Value = "[info] combat.round: hit chance calculated 81"
$ws2.Range("B5").Value = "[info] combat.armor_class: defender armor class is: 2"
$ws2.Range("B6").Value = "[info] combat.attack: inflicted 10"
$ws2.Range("B7").Value = "[info] combat.round: hit chance calculated 77"
$ws2.Range("B8").Value = "[info] combat.armor_class: defender armor class is: 6"
$ws2.Range("B9").Value = "[info] combat.attack: inflicted 2"

# Combat stat summary
$ws2.Range("B12").Value = "to_hit"
$ws2.Range("C12").Value = 81
$ws2.Range("B13").Value = "ac"
$ws2.Range("C13").Value = 2
$ws2.Range("B14").Value = "str"
$ws2.Range("C14").Value = 10
$ws2.Range("B15").Value = "def"
$ws2.Range("C15").Value = 3

# Make the new sheet the active one, with B16 selected (matches the saved view)
$ws2.Activate() | Out-Null
$ws2.Range("B16").Select() | Out-Null
